# 17.1.2 - add a new "2023" column (T) to the right of the existing
# 2007-2022 series, carrying over the header/value formatting from the
# last existing column (S), and tighten up the column widths so the new
# data columns read like the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new column T: header (row 4) + data value (row 5) -------------------
# Clone the formatting of the adjoining column S (bold header w/ medium
# top/bottom border for T4, "0.0" numeric style w/ same border for T5)
# and then overwrite just the value, same as Excel's own "copy format,
# then type a new number" workflow.
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5").Value = 75.099999999999994

$excel.CutCopyMode = $false

# --- column widths ---------------------------------------------------------
# Narrow the three label columns slightly and give every year column
# (D:T, now including the new T) an explicit, consistent width.
$ws.Columns("A:C").ColumnWidth = 34.666666666666664
$ws.Columns("D:T").ColumnWidth = 8

# --- reset the view's active cell back to the top-left --------------------
$null = $ws.Range("A1").Select()
